$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.659.88'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.744.40'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.94%  '
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.609'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.111'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.92%  '
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.391'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.09%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.68'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -16.41%  '
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.231.17'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.644.58'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000154'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.754.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.26'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '359.82'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.85'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.03%  '
$ws.Range("E22").Value = '  +4.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.989'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.171'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.66'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0925'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.99'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.09'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.00%  '
$ws.Range("E31").Value = '  +3.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '169.43'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.91%  '
$ws.Range("E33").Value = '  +0.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '20.52'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.00%  '
$ws.Range("E35").Value = '  +4.08%  '
$ws.Range("E36").Value = '  +0.92%  '
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("E38").Value = '  -1.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.19'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.89%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '332.40'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.90%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.09'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.69%  '
$ws.Range("E42").Value = '  +0.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.92'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0596'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.85'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.639'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.72%  '
$ws.Range("E47").Value = '  -0.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '136.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.20%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("E50").Value = '  +0.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.06'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.97%  '
